$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K header value (year 2022), same style as J4
$ws.Range("K4").Value = 2022

# New column K data values, matching the corresponding row's trend
$ws.Range("K5").Value = 1.6
$ws.Range("K6").Value = 0.4
$ws.Range("K7").Value = 0.9
$ws.Range("K8").Value = 0.6
$ws.Range("K9").Value = 2.1
$ws.Range("K10").Value = 0.6
$ws.Range("K11").Value = 0.9
$ws.Range("K12").Value = 2.2999999999999998
$ws.Range("K13").Value = 4.3
$ws.Range("K14").Value = 0.3

# Copy formatting (styles) from column J to column K for rows 4-14
$ws.Range("J4:J14").Copy()
$ws.Range("K4:K14").PasteSpecial(-4122)  # xlPasteFormats

# Update the selected cell (cosmetic, matches author's last selection)
$ws.Range("L7").Select()
